# "Game Zone_Issue Log.xlsx" update
#
# The resolution note on the usd_price / "0 and blank transactions" row
# was reworded: the value was actually deleted (not left in place), so
# the note is corrected from "left as it - ..." to "deleted - ...".
#
# Also restore the sheet's last-used selection (C6) / scroll position
# that the author had when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "deleted - no way to infer, would need to check with stakeholders"

$ws.Activate()
$ws.Range("C6").Select()
